# Update screen 0 and 1 -> button, cnt_up building
# The Translation sheet's text table (Table8) had a duplicated 5-row
# "cnt_up" counter block (rows 117:121) that needs to be removed; the
# remaining rows shift up, and the alignment for the now-shifted row 120
# needs to change from Right to Center.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the duplicated counter block (5 rows) - everything below shifts up.
$ws.Range("B117:B121").EntireRow.Delete() | Out-Null

# After the shift, former row 125 data is now on row 120; its Alignment
# (column D) changes from Right to Center.
$ws.Range("D120").Value = "Center"
